# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# sheet with refreshed quote data.
#
# The Price column cells are stored as plain text in the workbook (not
# numbers), so we prefix the assigned value with a leading apostrophe to
# force Excel to keep it as text instead of re-interpreting it as a
# number (which would silently drop meaningful trailing zeros, e.g.
# "1.00" -> 1, or introduce floating point artifacts, e.g.
# "588.56" -> 588.55999999999995). We then reset the cell Style back to
# "Normal" so the forced-text (quote-prefix) formatting doesn't linger
# on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.381.34"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.15%  "
$ws.Range("D3").Value = "'3.060.37"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.63%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.33%  "
$ws.Range("D5").Value = "'588.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").Value = "'155.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.11%  "
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("E8").Value = "  +1.44%  "
$ws.Range("D9").Value = "'3.059.24"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.26%  "
$ws.Range("E10").Value = "  -3.52%  "
$ws.Range("E11").Value = "  -1.54%  "
$ws.Range("E12").Value = "  -1.18%  "
$ws.Range("D13").Value = "'37.06"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.82%  "
$ws.Range("E14").Value = "  -3.54%  "
$ws.Range("E15").Value = "  -1.78%  "
$ws.Range("D16").Value = "'3.566.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.70%  "
$ws.Range("D17").Value = "'63.421.94"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.75%  "
$ws.Range("E18").Value = "  -1.71%  "
$ws.Range("D19").Value = "'3.057.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.51%  "
$ws.Range("D20").Value = "'473.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.28%  "
$ws.Range("D21").Value = "'14.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("E22").Value = "  -3.50%  "
$ws.Range("D23").Value = "'7.51"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.60%  "
$ws.Range("D24").Value = "'2.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.57%  "
$ws.Range("D25").Value = "'80.67"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("D26").Value = "'12.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.12%  "
$ws.Range("E27").Value = "  +4.03%  "
$ws.Range("D28").Value = "'0.998"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.24%  "
$ws.Range("D29").Value = "'7.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.08%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.33%  "
$ws.Range("E31").Value = "  -1.83%  "
$ws.Range("E32").Value = "  -2.59%  "
$ws.Range("E33").Value = "  -1.99%  "
$ws.Range("D34").Value = "'27.09"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.95%  "
$ws.Range("D35").Value = "'0.0₃0824"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.30%  "
$ws.Range("E36").Value = "  -1.93%  "
$ws.Range("D37").Value = "'3.32"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.54%  "
$ws.Range("D38").Value = "'5.98"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.81%  "
$ws.Range("D39").Value = "'2.21"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.84%  "
$ws.Range("E40").Value = "  -0.86%  "
$ws.Range("D41").Value = "'50.69"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.96%  "
$ws.Range("D42").Value = "'442.48"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.44%  "
$ws.Range("E43").Value = "  -0.56%  "
$ws.Range("D44").Value = "'41.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.12%  "
$ws.Range("E45").Value = "  +3.33%  "
$ws.Range("E46").Value = "  -4.05%  "
$ws.Range("D47").Value = "'2.791.55"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.32%  "
$ws.Range("D48").Value = "'130.38"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.50%  "
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("D50").Value = "'25.02"
$ws.Range("D50").Style = "Normal"
$ws.Range("E51").Value = "  +0.53%  "
